# Natmi following Dr Hou advice
# Rewrites the Sending/Target cluster combinations (ECs/FAPs/sCs) for the
# Ucn2 -> Crhr2 ligand-receptor pair and fills in the full 3x3 cluster matrix
# (rows 2-10) with the recalculated NATMI specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ 'A'=1; 'B'=2; 'C'=3; 'D'=4; 'E'=5; 'F'=6; 'G'=7; 'H'=8; 'I'=9; 'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16; 'Q'=17; 'R'=18; 'S'=19; 'T'=20 }

$rowData = @{}
$rowData[2] = @{ A='ECs'; B='Ucn2'; C='Crhr2'; D='ECs'; E=1.0; F=0.3333333333333333; G=0.3323133333333333; H=0.99694; I=0.2302327407427936; J=0.2302327407427935; K=1.0; L=0.3333333333333333; M=0.265436; N=0.796308; O=0.1143970297008367; P=0.1143970297008367; Q=0.08820792194666667; R=0.7938712975200001; S=0.0263379416808584; T=0.0263379416808584 }
$rowData[3] = @{ A='ECs'; B='Ucn2'; C='Crhr2'; D='FAPs'; E=1.0; F=0.3333333333333333; G=0.3323133333333333; H=0.99694; I=0.2302327407427936; J=0.2302327407427935; K=1.0; L=0.3333333333333333; M=0.055009; N=0.165027; O=0.02370765912239986; P=0.02370765912239986; Q=0.01828022415333333; R=0.16452201738; S=0.005458279336346012; T=0.00545827933634601 }
$rowData[4] = @{ A='ECs'; B='Ucn2'; C='Crhr2'; D='sCs'; E=1.0; F=0.3333333333333333; G=0.3323133333333333; H=0.99694; I=0.2302327407427936; J=0.2302327407427935; K=3.0; L=1.0; M=1.99986; N=5.99958; O=0.8618953111767634; P=0.8618953111767634; Q=0.6645801428; R=5.9812212852; S=0.1984365197255891; T=0.1984365197255891 }
$rowData[5] = @{ A='FAPs'; B='Ucn2'; C='Crhr2'; D='ECs'; E=3.0; F=1.0; G=0.476438; H=1.429314; I=0.3300849395169671; J=0.3300849395169671; K=1.0; L=0.3333333333333333; M=0.265436; N=0.796308; O=0.1143970297008367; P=0.1143970297008367; Q=0.126463796968; R=1.138174172712; S=0.03776073662972139; T=0.03776073662972139 }
$rowData[6] = @{ A='FAPs'; B='Ucn2'; C='Crhr2'; D='FAPs'; E=3.0; F=1.0; G=0.476438; H=1.429314; I=0.3300849395169671; J=0.3300849395169671; K=1.0; L=0.3333333333333333; M=0.055009; N=0.165027; O=0.02370765912239986; P=0.02370765912239986; Q=0.026208377942; R=0.235875401478; S=0.007825541227506231; T=0.00782554122750623 }
$rowData[7] = @{ A='FAPs'; B='Ucn2'; C='Crhr2'; D='sCs'; E=3.0; F=1.0; G=0.476438; H=1.429314; I=0.3300849395169671; J=0.3300849395169671; K=3.0; L=1.0; M=1.99986; N=5.99958; O=0.8618953111767634; P=0.8618953111767634; Q=0.95280929868; R=8.575283688119999; S=0.2844986616597395; T=0.2844986616597395 }
$rowData[8] = @{ A='sCs'; B='Ucn2'; C='Crhr2'; D='ECs'; E=3.0; F=1.0; G=0.6346286666666666; H=1.903886; I=0.4396823197402394; J=0.4396823197402393; K=1.0; L=0.3333333333333333; M=0.265436; N=0.796308; O=0.1143970297008367; P=0.1143970297008367; Q=0.1684532947653333; R=1.516079652888; S=0.05029835139025696; T=0.05029835139025696 }
$rowData[9] = @{ A='sCs'; B='Ucn2'; C='Crhr2'; D='FAPs'; E=3.0; F=1.0; G=0.6346286666666666; H=1.903886; I=0.4396823197402394; J=0.4396823197402393; K=1.0; L=0.3333333333333333; M=0.055009; N=0.165027; O=0.02370765912239986; P=0.02370765912239986; Q=0.03491028832466667; R=0.314192594922; S=0.01042383855854762; T=0.01042383855854762 }
$rowData[10] = @{ A='sCs'; B='Ucn2'; C='Crhr2'; D='sCs'; E=3.0; F=1.0; G=0.6346286666666666; H=1.903886; I=0.4396823197402394; J=0.4396823197402393; K=3.0; L=1.0; M=1.99986; N=5.99958; O=0.8618953111767634; P=0.8618953111767634; Q=1.26916848532; R=11.42251636788; S=0.3789601297914348; T=0.3789601297914347 }

foreach ($rowNum in 2..10) {
  $rowVals = $rowData[$rowNum]
  foreach ($col in @('A','B','C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T')) {
    if ($rowVals.ContainsKey($col)) {
      $ws.Cells.Item($rowNum, $colIndex[$col]).Value = $rowVals[$col]
    }
  }
}
